$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.525.59"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.98"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.52"
$ws.Range("E5").Value = "  -1.66%  "

# Row 6
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.47"
$ws.Range("E8").Value = "  +2.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("E9").Value = "  +3.05%  "

# Row 10
$ws.Range("E10").Value = "  +6.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.056.38"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.04"
$ws.Range("E13").Value = "  -3.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.795.77"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.641"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.546.57"
$ws.Range("E16").Value = "  +1.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("E17").Value = "  +1.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.11"
$ws.Range("E18").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.63"
$ws.Range("E19").Value = "  -1.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0797"
$ws.Range("E20").Value = "  +7.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  +5.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.73"
$ws.Range("E25").Value = "  +3.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.39"
$ws.Range("E26").Value = "  -1.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("B30").Value = "Swop.fi"
$ws.Range("C30").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "544.25"
$ws.Range("E30").Value = "  +941.12%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0525"
$ws.Range("E31").Value = "  +1.46%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.78"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  -0.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.60"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.89"
$ws.Range("E35").Value = "  +1.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.424.44"
$ws.Range("E36").Value = "  -2.16%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.644"
$ws.Range("E37").Value = "  +1.28%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.06"
$ws.Range("E38").Value = "  -0.49%  "

# Row 39
$ws.Range("E39").Value = "  +1.61%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.962"
$ws.Range("E40").Value = "  +6.42%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "82.69"
$ws.Range("E41").Value = "  -1.18%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").Value = "  -3.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$ws.Range("E44").Value = "  +2.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.05"
$ws.Range("E45").Value = "  +3.66%  "

# Row 46
$ws.Range("E46").Value = "  -1.02%  "

# Row 47
$ws.Range("E47").Value = "  -2.46%  "

# Row 48
$ws.Range("E48").Value = "  +3.14%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.949.96"
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.23"
$ws.Range("E50").Value = "  +6.98%  "

# Row 51
$ws.Range("E51").Value = "  -0.12%  "

